$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New timesheet entry: Saturday 2017-07-30 (Excel serial date 42946),
# 7 hours, with a note about building the first executable Jar,
# cleaning up colour coding, fixing a Java double-precision bug and
# publishing the first youtube video - describing the work that led to
# the properties-file-from-command-line feature mentioned in the commit.

# Set the values/formula for the new row first so that dependent
# formulas elsewhere on the sheet (e.g. F2's SUM(C:C) weekly total)
# recalculate against the final data.
$ws.Cells.Item(31, 1).Value = 42946
$ws.Cells.Item(31, 2).Formula = "=A31"
$ws.Cells.Item(31, 3).Value = 7
$ws.Cells.Item(31, 4).Value = "Built an executable Jar file, cleaned up the q state viewer colour coding, there may be a bug in how Java handles irrational negative doubles, learnt how to make a youtube video and published my first video on youtube."

# Copy the formatting (date format, wrap text, column styles, etc.) from
# the previous row (30) down onto the new row (31).
$ws.Range("A30:D30").Copy()
$ws.Range("A31:D31").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# The long note wraps onto 3 lines at this column width, same as other
# similarly long notes elsewhere in the sheet.
$ws.Rows.Item(31).RowHeight = 42.75

# Keep the same cell selected as before the edit.
$ws.Range("D31").Select()
